$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Reorder-RecordedBy($val) {
    $parts = $val -split ", "
    $idxSystem = [System.Array]::IndexOf($parts, "System")
    if ($idxSystem -lt 0) {
        return $val
    }

    # build $rest = all tokens except the (first) "System" token
    $rest = @()
    for ($i = 0; $i -lt $parts.Length; $i++) {
        if ($i -ne $idxSystem) {
            $rest += $parts[$i]
        }
    }

    # place "System" right after a lowercase "system" token if present,
    # otherwise put it at the very front
    $idxLower = [System.Array]::IndexOf($rest, "system")

    $newParts = @()
    if ($idxLower -ge 0) {
        for ($i = 0; $i -le $idxLower; $i++) { $newParts += $rest[$i] }
        $newParts += "System"
        for ($i = $idxLower + 1; $i -lt $rest.Length; $i++) { $newParts += $rest[$i] }
    } else {
        $newParts += "System"
        $newParts += $rest
    }

    return [string]::Join(", ", $newParts)
}

$lastRow = $ws.UsedRange.Rows.Count

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)
    $val = $cell.Value2
    if ($null -eq $val) { continue }

    $newVal = Reorder-RecordedBy $val
    if ($newVal -ne $val) {
        $cell.Value2 = $newVal
    }
}
